# Update recomputed TPM-derived NATMI ligand-receptor metrics
# (Tgm2-Itgb1) for rows 2-26, columns G-T, per commit "update scripts wuth new tpm".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 42.32036466666667
$ws.Range("H2").Value = 126.961094
$ws.Range("I2").Value = 0.285778576657872
$ws.Range("J2").Value = 0.2880046678857171
$ws.Range("M2").Value = 121.928739
$ws.Range("N2").Value = 365.786217
$ws.Range("O2").Value = 0.2282232151508951
$ws.Range("P2").Value = 0.2419720431319445
$ws.Range("Q2").Value = 5160.068697826822
$ws.Range("R2").Value = 46440.61828044139
$ws.Range("S2").Value = 0.06522130558610609
$ws.Range("T2").Value = 0.06968907791984411

# Row 3
$ws.Range("G3").Value = 42.32036466666667
$ws.Range("H3").Value = 126.961094
$ws.Range("I3").Value = 0.285778576657872
$ws.Range("J3").Value = 0.2880046678857171
$ws.Range("O3").Value = 0.2768624053389947
$ws.Range("P3").Value = 0.2935413991166814
$ws.Range("Q3").Value = 6259.788385025673
$ws.Range("R3").Value = 56338.09546523105
$ws.Range("S3").Value = 0.07912134412785274
$ws.Range("T3").Value = 0.08454129316330856

# Row 4
$ws.Range("G4").Value = 42.32036466666667
$ws.Range("H4").Value = 126.961094
$ws.Range("I4").Value = 0.285778576657872
$ws.Range("J4").Value = 0.2880046678857171
$ws.Range("M4").Value = 83.50496933333334
$ws.Range("N4").Value = 250.514908
$ws.Range("O4").Value = 0.1563025480180701
$ws.Range("P4").Value = 0.1657186665504434
$ws.Range("Q4").Value = 3533.960753665484
$ws.Range("R4").Value = 31805.64678298935
$ws.Range("S4").Value = 0.04466791970060276
$ws.Range("T4").Value = 0.04772774952232434

# Row 5
$ws.Range("G5").Value = 42.32036466666667
$ws.Range("H5").Value = 126.961094
$ws.Range("I5").Value = 0.285778576657872
$ws.Range("J5").Value = 0.2880046678857171
$ws.Range("M5").Value = 91.06846250000001
$ws.Range("N5").Value = 182.136925
$ws.Range("O5").Value = 0.1704597085236707
$ws.Range("P5").Value = 0.1204857969594293
$ws.Range("Q5").Value = 3854.050542632659
$ws.Range("R5").Value = 23124.30325579595
$ws.Range("S5").Value = 0.04871373287941036
$ws.Range("T5").Value = 0.03470047193824639

# Row 6
$ws.Range("G6").Value = 42.32036466666667
$ws.Range("H6").Value = 126.961094
$ws.Range("I6").Value = 0.285778576657872
$ws.Range("J6").Value = 0.2880046678857171
$ws.Range("M6").Value = 89.83562999999999
$ws.Range("N6").Value = 269.50689
$ws.Range("O6").Value = 0.1681521229683693
$ws.Range("P6").Value = 0.1782820942415013
$ws.Range("Q6").Value = 3801.87662165974
$ws.Range("R6").Value = 34216.88959493766
$ws.Range("S6").Value = 0.04805427436390005
$ws.Range("T6").Value = 0.0513460753419937

# Row 7
$ws.Range("I7").Value = 0.04213668412459876
$ws.Range("J7").Value = 0.04246491062777905
$ws.Range("M7").Value = 121.928739
$ws.Range("N7").Value = 365.786217
$ws.Range("O7").Value = 0.2282232151508951
$ws.Range("P7").Value = 0.2419720431319445
$ws.Range("Q7").Value = 760.8274466348779
$ws.Range("R7").Value = 6847.447019713901
$ws.Range("S7").Value = 0.009616569526713609
$ws.Range("T7").Value = 0.01027532118601912

# Row 8
$ws.Range("I8").Value = 0.04213668412459876
$ws.Range("J8").Value = 0.04246491062777905
$ws.Range("O8").Value = 0.2768624053389947
$ws.Range("P8").Value = 0.2935413991166814
$ws.Range("Q8").Value = 922.9758540733264
$ws.Range("S8").Value = 0.01166606371974585
$ws.Range("T8").Value = 0.0124652092790431

# Row 9
$ws.Range("I9").Value = 0.04213668412459876
$ws.Range("J9").Value = 0.04246491062777905
$ws.Range("M9").Value = 83.50496933333334
$ws.Range("N9").Value = 250.514908
$ws.Range("O9").Value = 0.1563025480180701
$ws.Range("P9").Value = 0.1657186665504434
$ws.Range("Q9").Value = 521.0656086519831
$ws.Range("R9").Value = 4689.590477867848
$ws.Range("S9").Value = 0.006586071093707348
$ws.Range("T9").Value = 0.007037228364419295

# Row 10
$ws.Range("I10").Value = 0.04213668412459876
$ws.Range("J10").Value = 0.04246491062777905
$ws.Range("M10").Value = 91.06846250000001
$ws.Range("N10").Value = 182.136925
$ws.Range("O10").Value = 0.1704597085236707
$ws.Range("P10").Value = 0.1204857969594293
$ws.Range("Q10").Value = 568.2613169060917
$ws.Range("R10").Value = 3409.56790143655
$ws.Range("S10").Value = 0.007182606894033089
$ws.Range("T10").Value = 0.005116418599798899

# Row 11
$ws.Range("I11").Value = 0.04213668412459876
$ws.Range("J11").Value = 0.04246491062777905
$ws.Range("M11").Value = 89.83562999999999
$ws.Range("N11").Value = 269.50689
$ws.Range("O11").Value = 0.1681521229683693
$ws.Range("P11").Value = 0.1782820942415013
$ws.Range("Q11").Value = 560.5685218292599
$ws.Range("R11").Value = 5045.116696463339
$ws.Range("S11").Value = 0.007085372890398865
$ws.Range("T11").Value = 0.007570733198498634

# Row 12
$ws.Range("G12").Value = 42.241047
$ws.Range("H12").Value = 126.723141
$ws.Range("I12").Value = 0.2852429647825406
$ws.Range("J12").Value = 0.2874648838260633
$ws.Range("M12").Value = 121.928739
$ws.Range("N12").Value = 365.786217
$ws.Range("O12").Value = 0.2282232151508951
$ws.Range("P12").Value = 0.2419720431319445
$ws.Range("Q12").Value = 5150.397594749733
$ws.Range("R12").Value = 46353.57835274759
$ws.Range("S12").Value = 0.06509906652184495
$ws.Range("T12").Value = 0.06955846526807961

# Row 13
$ws.Range("G13").Value = 42.241047
$ws.Range("H13").Value = 126.723141
$ws.Range("I13").Value = 0.2852429647825406
$ws.Range("J13").Value = 0.2874648838260633
$ws.Range("O13").Value = 0.2768624053389947
$ws.Range("P13").Value = 0.2935413991166814
$ws.Range("Q13").Value = 6248.056165503509
$ws.Range("R13").Value = 56232.50548953158
$ws.Range("S13").Value = 0.07897305333572034
$ws.Range("T13").Value = 0.0843828441952169

# Row 14
$ws.Range("G14").Value = 42.241047
$ws.Range("H14").Value = 126.723141
$ws.Range("I14").Value = 0.2852429647825406
$ws.Range("J14").Value = 0.2874648838260633
$ws.Range("M14").Value = 83.50496933333334
$ws.Range("N14").Value = 250.514908
$ws.Range("O14").Value = 0.1563025480180701
$ws.Range("P14").Value = 0.1657186665504434
$ws.Range("Q14").Value = 3527.337334342892
$ws.Range("R14").Value = 31746.03600908603
$ws.Range("S14").Value = 0.04458420219973971
$ws.Range("T14").Value = 0.04763829722773332

# Row 15
$ws.Range("G15").Value = 42.241047
$ws.Range("H15").Value = 126.723141
$ws.Range("I15").Value = 0.2852429647825406
$ws.Range("J15").Value = 0.2874648838260633
$ws.Range("M15").Value = 91.06846250000001
$ws.Range("N15").Value = 182.136925
$ws.Range("O15").Value = 0.1704597085236707
$ws.Range("P15").Value = 0.1204857969594293
$ws.Range("Q15").Value = 3846.827204680238
$ws.Range("R15").Value = 23080.96322808143
$ws.Range("S15").Value = 0.04862243263525954
$ws.Range("T15").Value = 0.034635435625633

# Row 16
$ws.Range("G16").Value = 42.241047
$ws.Range("H16").Value = 126.723141
$ws.Range("I16").Value = 0.2852429647825406
$ws.Range("J16").Value = 0.2874648838260633
$ws.Range("M16").Value = 89.83562999999999
$ws.Range("N16").Value = 269.50689
$ws.Range("O16").Value = 0.1681521229683693
$ws.Range("P16").Value = 0.1782820942415013
$ws.Range("Q16").Value = 3794.75106910461
$ws.Range("R16").Value = 34152.75962194149
$ws.Range("S16").Value = 0.04796421008997599
$ws.Range("T16").Value = 0.05124984150940044

# Row 17
$ws.Range("G17").Value = 3.4338745
$ws.Range("H17").Value = 6.867749
$ws.Range("I17").Value = 0.02318807445921414
$ws.Range("J17").Value = 0.0155791330048516
$ws.Range("M17").Value = 121.928739
$ws.Range("N17").Value = 365.786217
$ws.Range("O17").Value = 0.2282232151508951
$ws.Range("P17").Value = 0.2419720431319445
$ws.Range("Q17").Value = 418.6879876692554
$ws.Range("R17").Value = 2512.127926015533
$ws.Range("S17").Value = 0.005292056906240203
$ws.Range("T17").Value = 0.003769714643408251

# Row 18
$ws.Range("G18").Value = 3.4338745
$ws.Range("H18").Value = 6.867749
$ws.Range("I18").Value = 0.02318807445921414
$ws.Range("J18").Value = 0.0155791330048516
$ws.Range("O18").Value = 0.2768624053389947
$ws.Range("P18").Value = 0.2935413991166814
$ws.Range("Q18").Value = 507.9192459715849
$ws.Range("R18").Value = 3047.51547582951
$ws.Range("S18").Value = 0.006419906069957735
$ws.Range("T18").Value = 0.004573120499269006

# Row 19
$ws.Range("G19").Value = 3.4338745
$ws.Range("H19").Value = 6.867749
$ws.Range("I19").Value = 0.02318807445921414
$ws.Range("J19").Value = 0.0155791330048516
$ws.Range("M19").Value = 83.50496933333334
$ws.Range("N19").Value = 250.514908
$ws.Range("O19").Value = 0.1563025480180701
$ws.Range("P19").Value = 0.1657186665504434
$ws.Range("Q19").Value = 286.7455848170154
$ws.Range("R19").Value = 1720.473508902092
$ws.Range("S19").Value = 0.003624355121607902
$ws.Range("T19").Value = 0.002581753147576008

# Row 20
$ws.Range("G20").Value = 3.4338745
$ws.Range("H20").Value = 6.867749
$ws.Range("I20").Value = 0.02318807445921414
$ws.Range("J20").Value = 0.0155791330048516
$ws.Range("M20").Value = 91.06846250000001
$ws.Range("N20").Value = 182.136925
$ws.Range("O20").Value = 0.1704597085236707
$ws.Range("P20").Value = 0.1204857969594293
$ws.Range("Q20").Value = 312.7176711329563
$ws.Range("R20").Value = 1250.870684531825
$ws.Range("S20").Value = 0.003952632413542816
$ws.Range("T20").Value = 0.001877064256026494

# Row 21
$ws.Range("G21").Value = 3.4338745
$ws.Range("H21").Value = 6.867749
$ws.Range("I21").Value = 0.02318807445921414
$ws.Range("J21").Value = 0.0155791330048516
$ws.Range("M21").Value = 89.83562999999999
$ws.Range("N21").Value = 269.50689
$ws.Range("O21").Value = 0.1681521229683693
$ws.Range("P21").Value = 0.1782820942415013
$ws.Range("Q21").Value = 308.484279048435
$ws.Range("R21").Value = 1850.90567429061
$ws.Range("S21").Value = 0.003899123947865479
$ws.Range("T21").Value = 0.002777480458571835

# Row 22
$ws.Range("G22").Value = 53.85273233333334
$ws.Range("H22").Value = 161.558197
$ws.Range("I22").Value = 0.3636536999757743
$ws.Range("J22").Value = 0.3664864046555889
$ws.Range("M22").Value = 121.928739
$ws.Range("N22").Value = 365.786217
$ws.Range("O22").Value = 0.2282232151508951
$ws.Range("P22").Value = 0.2419720431319445
$ws.Range("Q22").Value = 6566.195745107861
$ws.Range("R22").Value = 59095.76170597075
$ws.Range("S22").Value = 0.08299421660999019
$ws.Range("T22").Value = 0.08867946411459343

# Row 23
$ws.Range("G23").Value = 53.85273233333334
$ws.Range("H23").Value = 161.558197
$ws.Range("I23").Value = 0.3636536999757743
$ws.Range("J23").Value = 0.3664864046555889
$ws.Range("O23").Value = 0.2768624053389947
$ws.Range("P23").Value = 0.2935413991166814
$ws.Range("Q23").Value = 7965.590821754336
$ws.Range("R23").Value = 71690.31739578902
$ws.Range("S23").Value = 0.100682038085718
$ws.Range("T23").Value = 0.1075789319798438

# Row 24
$ws.Range("G24").Value = 53.85273233333334
$ws.Range("H24").Value = 161.558197
$ws.Range("I24").Value = 0.3636536999757743
$ws.Range("J24").Value = 0.3664864046555889
$ws.Range("M24").Value = 83.50496933333334
$ws.Range("N24").Value = 250.514908
$ws.Range("O24").Value = 0.1563025480180701
$ws.Range("P24").Value = 0.1657186665504434
$ws.Range("Q24").Value = 4496.970762011209
$ws.Range("R24").Value = 40472.73685810088
$ws.Range("S24").Value = 0.05683999990241231
$ws.Range("T24").Value = 0.06073363828839039

# Row 25
$ws.Range("G25").Value = 53.85273233333334
$ws.Range("H25").Value = 161.558197
$ws.Range("I25").Value = 0.3636536999757743
$ws.Range("J25").Value = 0.3664864046555889
$ws.Range("M25").Value = 91.06846250000001
$ws.Range("N25").Value = 182.136925
$ws.Range("O25").Value = 0.1704597085236707
$ws.Range("P25").Value = 0.1204857969594293
$ws.Range("Q25").Value = 4904.285535020705
$ws.Range("R25").Value = 29425.71321012423
$ws.Range("S25").Value = 0.0619883037014249
$ws.Range("T25").Value = 0.04415640653972455

# Row 26
$ws.Range("G26").Value = 53.85273233333334
$ws.Range("H26").Value = 161.558197
$ws.Range("I26").Value = 0.3636536999757743
$ws.Range("J26").Value = 0.3664864046555889
$ws.Range("M26").Value = 89.83562999999999
$ws.Range("N26").Value = 269.50689
$ws.Range("O26").Value = 0.1681521229683693
$ws.Range("P26").Value = 0.1782820942415013
$ws.Range("Q26").Value = 4837.89413638637
$ws.Range("R26").Value = 43541.04722747733
$ws.Range("S26").Value = 0.06114914167622888
$ws.Range("T26").Value = 0.06533796373303669

